# Update LR-pair TPM-derived stats on the "Ltb-Ltbr" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6577333333333334
$ws.Range("H2").Value = 1.9732
$ws.Range("M2").Value = 5.855966
$ws.Range("N2").Value = 17.567898
$ws.Range("O2").Value = 0.2049970778295038
$ws.Range("P2").Value = 0.2049970778295038
$ws.Range("Q2").Value = 3.851664037066667
$ws.Range("R2").Value = 34.6649763336
$ws.Range("S2").Value = 0.2049970778295038
$ws.Range("T2").Value = 0.2049970778295038

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6577333333333334
$ws.Range("H3").Value = 1.9732
$ws.Range("O3").Value = 0.5343457573275825
$ws.Range("P3").Value = 0.5343457573275825
$ws.Range("Q3").Value = 10.03975451088889
$ws.Range("R3").Value = 90.35779059800001
$ws.Range("S3").Value = 0.5343457573275825
$ws.Range("T3").Value = 0.5343457573275825

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.6577333333333334
$ws.Range("H4").Value = 1.9732
$ws.Range("M4").Value = 7.445957333333333
$ws.Range("N4").Value = 22.337872
$ws.Range("O4").Value = 0.2606571648429137
$ws.Range("P4").Value = 0.2606571648429137
$ws.Range("Q4").Value = 4.897454336711111
$ws.Range("R4").Value = 44.0770890304
$ws.Range("S4").Value = 0.2606571648429137
$ws.Range("T4").Value = 0.2606571648429137
